# Apply odds updates to "Jogos da Semana" worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("H5").Value = 3.3
$ws.Range("AH5").Value = 21
$ws.Range("AP5").Value = 26

# Row 7
$ws.Range("G7").Value = 1.8
$ws.Range("H7").Value = 3.4
$ws.Range("I7").Value = 4.75
$ws.Range("N7").Value = 9
$ws.Range("Q7").Value = 2.05
$ws.Range("R7").Value = 1.75
$ws.Range("Y7").Value = 8.5
$ws.Range("AC7").Value = 9
$ws.Range("AG7").Value = 12
$ws.Range("AP7").Value = 21
$ws.Range("AV7").Value = 51

# Row 10
$ws.Range("G10").Value = 2.5
$ws.Range("H10").Value = 3.1
$ws.Range("I10").Value = 3.1
$ws.Range("Z10").Value = 23
$ws.Range("AB10").Value = 34
$ws.Range("AF10").Value = 51
$ws.Range("AG10").Value = 8
$ws.Range("AK10").Value = 26
$ws.Range("AM10").Value = 401
$ws.Range("AU10").Value = 8.5
$ws.Range("AX10").Value = 17
$ws.Range("AZ10").Value = 51
$ws.Range("BB10").Value = 251

# Row 14
$ws.Range("M14").Value = 1.06
$ws.Range("N14").Value = 10
$ws.Range("Q14").Value = 1.95
$ws.Range("R14").Value = 1.95

# Row 20
$ws.Range("H20").Value = 7.9
$ws.Range("K20").Value = 3.65
$ws.Range("L20").Value = 1.27
$ws.Range("P20").Value = 7
$ws.Range("Q20").Value = 1.23
$ws.Range("R20").Value = 3.7
$ws.Range("S20").Value = 1.14
$ws.Range("T20").Value = 4.85
$ws.Range("U20").Value = 2.22
$ws.Range("V20").Value = 1.6
$ws.Range("AC20").Value = 27
$ws.Range("AG20").Value = 13
$ws.Range("AH20").Value = 7.8
$ws.Range("AJ20").Value = 6.6
$ws.Range("AL20").Value = 37
$ws.Range("AT20").Value = 4.85
$ws.Range("AV20").Value = 80
$ws.Range("AW20").Value = 3.5
$ws.Range("BA20").Value = 22
